$wb = $excel.ActiveWorkbook

# The "settings" sheet's A2 cell held the old form title string; the author
# repurposed that string for something else ("Take a chance?") while the
# "survey" sheet's C2 cell keeps its original "Your email" text.
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("A2").Value = "Take a chance?"

# Make "settings" the active/selected tab (was "survey"), with A2 selected.
$wsSettings.Activate()
$wsSettings.Range("A2").Select()
